# 2023-10-30: added PRL with Abou
#
# The "Random traction yielding transition in epithelial tissues" preprint
# (previously row 17) was published in Phys. Rev. Lett.  Its row is updated
# with the publication details and moved so the table stays ordered by
# publication date: the old preprint row is removed and a fresh row is
# inserted just above the "Emergent chirality..." row, which shifts that
# row (and the rows below it) back down to where they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "preprint" row for the Amiri et al. paper (row 17) ...
$ws.Rows("17:17").Delete() | Out-Null

# ... and insert a fresh row above the "Emergent chirality..." row (row 16)
# to hold the now-published version, restoring the table to 20 rows with
# "Emergent chirality..." back at row 17.
$ws.Rows("16:16").Insert() | Out-Null

$ws.Range("A16").Value = "2023-10-30"
$ws.Range("B16").Value = "Random traction yielding transition in epithelial tissues"
$ws.Range("C16").Value = "A. Amiri, C. Duclut, F. Jülicher, M. Popović"
$ws.Range("D16").Value = "Phys. Rev. Lett."
$ws.Range("E16").Value = 131
$ws.Range("F16").Value = 188401
$ws.Range("G16").Value = "We investigate how randomly oriented cell traction forces lead to fluidisation in a vertex model of epithelial tissues. We find that the fluidisation occurs at a critical value of the traction force magnitude `$F_c`$. We show that this transition exhibits critical behaviour, similar to the yielding transition of sheared amorphous solids. However, we find that it belongs to a different universality class, even though it satisfies the same scaling relations between critical exponents established in the yielding transition of sheared amorphous solids. Our work provides a fluidisation mechanism through active force generation that could be relevant in biological tissues."
$ws.Range("H16").Value = "amiri2023random"
$ws.Range("I16").Value = "https://link.aps.org/doi/10.1103/PhysRevLett.131.188401"
$ws.Range("J16").Value = "2211.02159"
$ws.Range("K16").Value = "rheology, tissues, mechanics, vertexModel, yielding"

# Column C got a bit narrower, and the new volume/page columns need room,
# so column D grows a dedicated width.
$ws.Columns("C").ColumnWidth = 70.6640625
$ws.Columns("D").ColumnWidth = 14.6640625

# Reflect the author's last selection when the file was saved.
$ws.Range("C22").Select() | Out-Null
